$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextSafe($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "62.389.38"
$ws.Range("E2").Value = "  -1.28%  "

$ws.Range("D3").Value = "2.429.53"
$ws.Range("E3").Value = "  -0.98%  "

$ws.Range("E4").Value = "  -0.29%  "

Set-TextSafe "D5" "570.67"
$ws.Range("E5").Value = "  -0.15%  "

Set-TextSafe "D6" "143.07"
$ws.Range("E6").Value = "  -2.83%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  -1.45%  "

$ws.Range("D9").Value = "2.425.85"
$ws.Range("E9").Value = "  -1.42%  "

$ws.Range("E10").Value = "  -4.41%  "

$ws.Range("E11").Value = "  +0.80%  "

$ws.Range("E12").Value = "  -1.06%  "

$ws.Range("E13").Value = "  -1.90%  "

Set-TextSafe "D14" "26.49"
$ws.Range("E14").Value = "  -1.78%  "

$ws.Range("E15").Value = "  -4.37%  "

$ws.Range("D16").Value = "2.868.39"
$ws.Range("E16").Value = "  -1.76%  "

$ws.Range("D17").Value = "62.230.67"
$ws.Range("E17").Value = "  -1.63%  "

$ws.Range("D18").Value = "2.419.19"
$ws.Range("E18").Value = "  -1.63%  "

Set-TextSafe "D19" "11.01"
$ws.Range("E19").Value = "  -3.26%  "

Set-TextSafe "D20" "7.12"
$ws.Range("E20").Value = "  -2.54%  "

Set-TextSafe "D21" "324.21"
$ws.Range("E21").Value = "  -1.02%  "

$ws.Range("E22").Value = "  -2.00%  "

$ws.Range("E23").Value = "  +2.97%  "

$ws.Range("E24").Value = "  +0.73%  "

Set-TextSafe "D25" "65.11"
$ws.Range("E25").Value = "  -2.93%  "

Set-TextSafe "D26" "621.42"
$ws.Range("E26").Value = "  -0.57%  "

Set-TextSafe "D27" "8.99"
$ws.Range("E27").Value = "  +1.49%  "

$ws.Range("D28").Value = "0.0₃0959"
$ws.Range("E28").Value = "  -7.45%  "

$ws.Range("E29").Value = "  -1.31%  "

Set-TextSafe "D30" "1.00"
$ws.Range("E30").Value = "  +0.60%  "

$ws.Range("E31").Value = "  -4.03%  "

Set-TextSafe "D32" "7.98"
$ws.Range("E32").Value = "  -3.59%  "

$ws.Range("E33").Value = "  -2.35%  "

Set-TextSafe "D34" "0.134"
$ws.Range("E34").Value = "  -7.69%  "

Set-TextSafe "D35" "5.02"
$ws.Range("E35").Value = "  -2.23%  "

$ws.Range("E36").Value = "  +0.23%  "

$ws.Range("E37").Value = "  -4.90%  "

$ws.Range("E38").Value = "  -2.19%  "

$ws.Range("E39").Value = "  -1.54%  "

Set-TextSafe "D40" "147.07"
$ws.Range("E40").Value = "  +0.96%  "

$ws.Range("E41").Value = "  -4.13%  "

$ws.Range("E42").Value = "  -5.13%  "

Set-TextSafe "D43" "42.34"
$ws.Range("E43").Value = "  +1.01%  "

$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("E45").Value = "  -7.10%  "

Set-TextSafe "D46" "144.58"
$ws.Range("E46").Value = "  -2.87%  "

$ws.Range("E47").Value = "  -1.48%  "

$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextSafe "D48" "0.0521"
$ws.Range("E48").Value = "  -3.80%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextSafe "D49" "19.92"
$ws.Range("E49").Value = "  -4.24%  "

Set-TextSafe "D50" "0.593"
$ws.Range("E50").Value = "  -1.89%  "

$ws.Range("E51").Value = "  -3.50%  "
